# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect the latest generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of event name -> new value for column F ("想去人数").
# Using the event name (column C) as the key makes the update robust
# regardless of each sheet's exact row layout.
$updates = @{
    "宜春·第三十五届静卿国风动漫文化展览会"                                   = 624
    "万载·第八届馨缘动漫文化展"                                               = 75
    "余干·沧渊动漫游戏嘉年华"                                                 = 45
    "南昌·幻梦境国际动漫游戏嘉年华1th"                                        = 4573
    "吉安·COMIC LIFE周年庆典"                                                = 1825
    "宜春·星语国风动漫游戏博览会"                                             = 128
    "景德镇·第十五届瓷都ACG动漫游戏博览会"                                    = 3082
    "樟树·第二届静卿国风动漫文化展览会"                                       = 581
    "赣州·第一届环梦动漫游戏嘉年华"                                           = 602
    "上饶·第十五届IX Group国风嘉年华暨十周年庆典"                            = 517
    "九江·第一届异次元动漫嘉年华"                                             = 515
    "上饶·囧喵喵国风动漫展"                                                   = 367
    "南昌·第一届异次元动漫嘉年华"                                             = 1763
    "赣州·第二届异次元动漫嘉年华"                                             = 1309
    "信丰·七夕节UPUP动漫展"                                                  = 117
    "南昌·W·MEETING动漫游戏盛典"                                             = 1557
    "吉安·WF无线次元新星动漫博览会"                                           = 603
    "赣州·十万伏特-星铁&音乐 次元音乐同人only2.0"                            = 43
    "高安·第二届静卿国风动漫文化展览会"                                       = 526
    "九江·动漫畅想（取消）"                                                   = 24
    "萍乡·夏花Flower·2024夏季国漫展"                                         = 89
    "南昌·CM03·配音演员孙路路专场见面会"                                     = 82
    "南昌·CM03动漫游戏博览会"                                                = 3598
    "南昌·第四届龙年动漫展——暑假最后的狂欢"                                  = 747
    "赣州·第五人格only"                                                      = 64
    "南昌·Sunflower Garden动漫游戏展"                                        = 253
    "南昌·萌卡动漫展"                                                         = 1722
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Determine the last used row from column C (names).
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp = -4162

    for ($r = 2; $r -le $lastRow; $r++) {
        # NOTE: use Value2 (not the parameterized Value property) - reading
        # Value directly without an index returns its property descriptor
        # rather than invoking the getter in this COM shim.
        $name = $ws.Cells.Item($r, 3).Value2
        if ($null -ne $name -and $updates.ContainsKey($name)) {
            $ws.Cells.Item($r, 6).Value = $updates[$name]
        }
    }
}
